# Evaluation.xlsx -- redaction of company name + criteria table rework
# (Price row removed, Quality and OHS Policies row added, weightings and
#  comments updated, row heights/selection refreshed.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Criteria column (A) -- "Price" dropped, rows shift up, new row added
#    (A6 is set last so the new shared strings are appended in the same
#    order as the target workbook: comments first, then the new criteria)
# ---------------------------------------------------------------------
$ws.Range("A3").Value = "Technical Proposal"
$ws.Range("A4").Value = "Timeline"
$ws.Range("A5").Value = "References"

# ---------------------------------------------------------------------
# 2. Weighting column (B)
# ---------------------------------------------------------------------
$ws.Range("B2").Value = 35
$ws.Range("B3").Value = 30
$ws.Range("B4").Value = 15
$ws.Range("B5").Value = 10
$ws.Range("B6").Value = 10

# ---------------------------------------------------------------------
# 3. Comments column (D) -- new guidance text
# ---------------------------------------------------------------------
$ws.Range("D2").Value = "at least 10 years is acceptable"
$ws.Range("D5").Value = "a high score will be given for the offer to provide references upon contact, no comment is unacceptable"
$ws.Range("D6").Value = "alignment with a standard is expected"

# New criteria row (appends the last shared string, matching the order
# in which the target workbook's sharedStrings table was built)
$ws.Range("A6").Value = "Quality and OHS Policies"

# ---------------------------------------------------------------------
# 4. C6 is fully cleared (no cell left behind at all)
# ---------------------------------------------------------------------
$ws.Range("C6").Clear()

# ---------------------------------------------------------------------
# 5. Row heights
# ---------------------------------------------------------------------
$ws.Rows.Item(1).AutoFit()
$ws.Rows.Item(2).RowHeight = 28.8
$ws.Rows.Item(3).RowHeight = 28.8
$ws.Rows.Item(4).AutoFit()
$ws.Rows.Item(5).RowHeight = 100.8
$ws.Rows.Item(6).RowHeight = 43.2

# ---------------------------------------------------------------------
# 6. Selection moves to A7
# ---------------------------------------------------------------------
$ws.Range("A7").Select()
